# Semana 32 de 2025 - update poisson.xlsx data (Esperado/Observado/valor p)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 0.37

$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.37

$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 0

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 0.18

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 11

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0.37

$ws.Range("C9").Value = 43
$ws.Range("D9").Value = 45
$ws.Range("E9").Value = 0.06

$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 0.04

$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.05

$ws.Range("C14").Value = 7

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0.37

$ws.Range("C16").Value = 0
$ws.Range("E16").Value = 1

$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 1

$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 15
$ws.Range("E18").Value = 0

$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 0.14

$ws.Range("D22").Value = 0

$ws.Range("D23").Value = 4

$ws.Range("D25").Value = 4

$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0.37

$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0.27

$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 1

$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 4
$ws.Range("E31").Value = 0.09

$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 0.01

$ws.Range("C35").Value = 7
$ws.Range("D35").Value = 5

$wb.Save()
